# LOM3220.xlsx update — build 2023-01-09
# 1) Ativação date bumped from 01/01/2012 -> 01/01/2023 (cells B8/C8/B13/C13)
# 2) English "Objectives:" text filled in (B11/C11)
# 3) English "Short syllabus:" text filled in (B14/C14)
# 4) English "Syllabus:" text filled in (B16/C16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues / xlPasteFormats paste-special codes
$xlPasteValues = -4163
$xlPasteFormats = -4122

# --- 1) Update the activation date in place (B8, C8, B13, C13 all share the
# same string). Assigning the literal text directly would make Excel's
# "looks like a date" auto-detection convert it into a date serial, so we
# route it through a text formula and collapse that formula down to its
# literal (text) result with a Paste Special > Values, which keeps the cell
# a plain text cell instead of a date.
foreach ($addr in @("B8", "C8", "B13", "C13")) {
    $cell = $ws.Range($addr)
    $cell.Formula = '="01/01/2023"'
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}
$excel.CutCopyMode = $false

# --- 2) Objectives (English) — row 11 previously only had column A populated
$ws.Range("B11").Value = "To present the concepts of spintronics and the potential applications in quantum computing."
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)

$ws.Range("C11").Value = "To present the concepts of spintronics and the potential applications in quantum computing."
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)

# --- 3) Short syllabus (English) — row 14 previously only had column A populated
$ws.Range("B14").Value = "Introduction to nanotechnology. Metal spintronics. Semiconductor spintronics. Spintronics devices. Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms. decoherence. Quantum Dots. Kane transistor. Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("B15").Copy()
$ws.Range("B14").PasteSpecial($xlPasteFormats)

$ws.Range("C14").Value = "Introduction to nanotechnology. Metal spintronics. Semiconductor spintronics. Spintronics devices. Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms. decoherence. Quantum Dots. Kane transistor. Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("C15").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)

# --- 4) Syllabus (English) — row 16 previously only had column A populated
$ws.Range("B16").Value = "Introduction to nanotechnology.Metal spintronics. Semiconductor SpintronicsSpintronics devices.Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms.decoherence. Quantum Dots.Kane transistor.Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)

$ws.Range("C16").Value = "Introduction to nanotechnology.Metal spintronics. Semiconductor SpintronicsSpintronics devices.Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms.decoherence. Quantum Dots.Kane transistor.Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
